# Estadisticos Matutinos 15 Oct
# Fill in the grade statistics (Aprobados/Reprobados/Por_Apro/Por_Repro/Promedio/
# Blancos/Por_Blan) for rows 8-12 (Medina Tolentino Elio / Polanco Dominguez Rosa
# Maria) on the "1er Parcial" and "3er Parcial" sheets. "2o Parcial" keeps its
# placeholder values.

$wb = $excel.ActiveWorkbook

# Columns: D=Totales(unchanged) E=Aprobados F=Reprobados G=Por_Apro H=Por_Repro
#          I=Promedio(new) J=Blancos K=Por_Blan
$data = @{
    8  = @{ E = 28; F = 10; G = 73.68; H = 26.32; I = 7;   J = 0;  K = 0 }
    9  = @{ E = 3;  F = 20; G = 13.04; H = 86.96; I = 8.7; J = 20; K = 86.96 }
    10 = @{ E = 3;  F = 28; G = 9.68;  H = 90.32; I = 9.3; J = 28; K = 90.32 }
    11 = @{ E = 6;  F = 32; G = 15.79; H = 84.21; I = 9.2; J = 32; K = 84.21 }
    12 = @{ E = 3;  F = 30; G = 9.09;  H = 90.91; I = 9;   J = 30; K = 90.91 }
}

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $data.Keys) {
        $vals = $data[$row]
        foreach ($col in $vals.Keys) {
            $ws.Range("$col$row").Value = $vals[$col]
        }
    }
}
